$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 109 values (weekly insert pushes existing rows down)
$ws.Range("D109").Value = 44957
$ws.Range("J109").Value = 2000
$ws.Range("K109").Value = 3000
$ws.Range("L109").Value = 3500
$ws.Range("M109").Value = 3250
$ws.Range("P109").Value = 2167

# Cascade shift of pre-existing rows down by one position (109-198 -> 110-199)
$ws.Range("D110").Value = 44355
$ws.Range("J110").Value = 3200
$ws.Range("K110").Value = 1300
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 1400
$ws.Range("P110").Value = 933
$ws.Range("D111").Value = 44488
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 1300
$ws.Range("L111").Value = 1500
$ws.Range("M111").Value = 1400
$ws.Range("P111").Value = 933
$ws.Range("D112").Value = 44299
$ws.Range("J112").Value = 3400
$ws.Range("K112").Value = 2000
$ws.Range("L112").Value = 2500
$ws.Range("M112").Value = 2250
$ws.Range("P112").Value = 1500
$ws.Range("D113").Value = 44915
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 2750
$ws.Range("P113").Value = 1833
$ws.Range("D114").Value = 44642
$ws.Range("J114").Value = 2400
$ws.Range("K114").Value = 2300
$ws.Range("L114").Value = 2500
$ws.Range("M114").Value = 2400
$ws.Range("P114").Value = 1600
$ws.Range("D115").Value = 44649
$ws.Range("J115").Value = 2800
$ws.Range("K115").Value = 2300
$ws.Range("L115").Value = 2500
$ws.Range("M115").Value = 2400
$ws.Range("P115").Value = 1600
$ws.Range("D116").Value = 44588
$ws.Range("J116").Value = 3200
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 2750
$ws.Range("P116").Value = 1833
$ws.Range("D117").Value = 44278
$ws.Range("J117").Value = 3400
$ws.Range("K117").Value = 2000
$ws.Range("L117").Value = 2500
$ws.Range("M117").Value = 2250
$ws.Range("P117").Value = 1500
$ws.Range("D118").Value = 44292
$ws.Range("J118").Value = 3400
$ws.Range("K118").Value = 2000
$ws.Range("L118").Value = 2500
$ws.Range("M118").Value = 2250
$ws.Range("P118").Value = 1500
$ws.Range("D119").Value = 44392
$ws.Range("J119").Value = 3320
$ws.Range("K119").Value = 1500
$ws.Range("L119").Value = 2000
$ws.Range("M119").Value = 1750
$ws.Range("P119").Value = 1167
$ws.Range("D120").Value = 44308
$ws.Range("J120").Value = 3200
$ws.Range("K120").Value = 1300
$ws.Range("L120").Value = 1500
$ws.Range("M120").Value = 1400
$ws.Range("P120").Value = 933
$ws.Range("D121").Value = 44614
$ws.Range("J121").Value = 2300
$ws.Range("K121").Value = 2300
$ws.Range("L121").Value = 2500
$ws.Range("M121").Value = 2400
$ws.Range("P121").Value = 1600
$ws.Range("D122").Value = 44747
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 1500
$ws.Range("L122").Value = 2000
$ws.Range("M122").Value = 1750
$ws.Range("P122").Value = 1167
$ws.Range("D123").Value = 44845
$ws.Range("J123").Value = 2800
$ws.Range("K123").Value = 2000
$ws.Range("L123").Value = 2500
$ws.Range("M123").Value = 2250
$ws.Range("P123").Value = 1500
$ws.Range("D124").Value = 44434
$ws.Range("J124").Value = 3360
$ws.Range("K124").Value = 2000
$ws.Range("L124").Value = 2500
$ws.Range("M124").Value = 2250
$ws.Range("P124").Value = 1500
$ws.Range("D125").Value = 44483
$ws.Range("J125").Value = 3300
$ws.Range("K125").Value = 1500
$ws.Range("L125").Value = 2000
$ws.Range("M125").Value = 1750
$ws.Range("P125").Value = 1167
$ws.Range("D126").Value = 44817
$ws.Range("J126").Value = 2600
$ws.Range("K126").Value = 2000
$ws.Range("L126").Value = 2500
$ws.Range("M126").Value = 2250
$ws.Range("P126").Value = 1500
$ws.Range("D127").Value = 44931
$ws.Range("J127").Value = 2400
$ws.Range("K127").Value = 3500
$ws.Range("L127").Value = 4000
$ws.Range("M127").Value = 3750
$ws.Range("P127").Value = 2500
$ws.Range("D128").Value = 44168
$ws.Range("J128").Value = 2800
$ws.Range("K128").Value = 1300
$ws.Range("L128").Value = 1500
$ws.Range("M128").Value = 1400
$ws.Range("P128").Value = 933
$ws.Range("D129").Value = 44245
$ws.Range("J129").Value = 3200
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 2000
$ws.Range("M129").Value = 1750
$ws.Range("P129").Value = 1167
$ws.Range("D130").Value = 44455
$ws.Range("J130").Value = 3200
$ws.Range("K130").Value = 2000
$ws.Range("L130").Value = 2500
$ws.Range("M130").Value = 2250
$ws.Range("P130").Value = 1500
$ws.Range("D131").Value = 44719
$ws.Range("J131").Value = 3200
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2000
$ws.Range("M131").Value = 1750
$ws.Range("P131").Value = 1167
$ws.Range("D132").Value = 44803
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 2000
$ws.Range("L132").Value = 2500
$ws.Range("M132").Value = 2250
$ws.Range("P132").Value = 1500
$ws.Range("D133").Value = 44315
$ws.Range("J133").Value = 3120
$ws.Range("K133").Value = 1300
$ws.Range("L133").Value = 1500
$ws.Range("M133").Value = 1400
$ws.Range("P133").Value = 933
$ws.Range("D134").Value = 44490
$ws.Range("J134").Value = 3200
$ws.Range("K134").Value = 1300
$ws.Range("L134").Value = 1500
$ws.Range("M134").Value = 1400
$ws.Range("P134").Value = 933
$ws.Range("D135").Value = 44656
$ws.Range("J135").Value = 2400
$ws.Range("K135").Value = 2000
$ws.Range("L135").Value = 2500
$ws.Range("M135").Value = 2250
$ws.Range("P135").Value = 1500
$ws.Range("D136").Value = 44413
$ws.Range("J136").Value = 3360
$ws.Range("K136").Value = 2000
$ws.Range("L136").Value = 2500
$ws.Range("M136").Value = 2250
$ws.Range("P136").Value = 1500
$ws.Range("D137").Value = 44894
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 1500
$ws.Range("L137").Value = 2000
$ws.Range("M137").Value = 1750
$ws.Range("P137").Value = 1167
$ws.Range("D138").Value = 44306
$ws.Range("J138").Value = 3400
$ws.Range("K138").Value = 2000
$ws.Range("L138").Value = 2500
$ws.Range("M138").Value = 2250
$ws.Range("P138").Value = 1500
$ws.Range("D139").Value = 44462
$ws.Range("J139").Value = 3200
$ws.Range("K139").Value = 1500
$ws.Range("L139").Value = 2000
$ws.Range("M139").Value = 1750
$ws.Range("P139").Value = 1167
$ws.Range("D140").Value = 44938
$ws.Range("J140").Value = 2500
$ws.Range("K140").Value = 3500
$ws.Range("L140").Value = 4000
$ws.Range("M140").Value = 3750
$ws.Range("P140").Value = 2500
$ws.Range("D141").Value = 44698
$ws.Range("J141").Value = 3200
$ws.Range("K141").Value = 1500
$ws.Range("L141").Value = 2000
$ws.Range("M141").Value = 1750
$ws.Range("P141").Value = 1167
$ws.Range("D142").Value = 44565
$ws.Range("J142").Value = 3000
$ws.Range("K142").Value = 3000
$ws.Range("L142").Value = 3500
$ws.Range("M142").Value = 3250
$ws.Range("P142").Value = 2167
$ws.Range("D143").Value = 44348
$ws.Range("J143").Value = 3360
$ws.Range("K143").Value = 1300
$ws.Range("L143").Value = 1500
$ws.Range("M143").Value = 1400
$ws.Range("P143").Value = 933
$ws.Range("D144").Value = 44343
$ws.Range("J144").Value = 3340
$ws.Range("K144").Value = 1300
$ws.Range("L144").Value = 1500
$ws.Range("M144").Value = 1400
$ws.Range("P144").Value = 933
$ws.Range("D145").Value = 44848
$ws.Range("J145").Value = 3000
$ws.Range("K145").Value = 1500
$ws.Range("L145").Value = 2000
$ws.Range("M145").Value = 1750
$ws.Range("P145").Value = 1167
$ws.Range("D146").Value = 44294
$ws.Range("J146").Value = 3000
$ws.Range("K146").Value = 2000
$ws.Range("L146").Value = 2500
$ws.Range("M146").Value = 2250
$ws.Range("P146").Value = 1500
$ws.Range("D147").Value = 44518
$ws.Range("J147").Value = 3320
$ws.Range("K147").Value = 1300
$ws.Range("L147").Value = 1500
$ws.Range("M147").Value = 1400
$ws.Range("P147").Value = 933
$ws.Range("D148").Value = 44726
$ws.Range("J148").Value = 3200
$ws.Range("K148").Value = 1500
$ws.Range("L148").Value = 2000
$ws.Range("M148").Value = 1750
$ws.Range("P148").Value = 1167
$ws.Range("D149").Value = 44194
$ws.Range("J149").Value = 3300
$ws.Range("K149").Value = 1400
$ws.Range("L149").Value = 1500
$ws.Range("M149").Value = 1450
$ws.Range("P149").Value = 967
$ws.Range("D150").Value = 44602
$ws.Range("J150").Value = 2960
$ws.Range("K150").Value = 2300
$ws.Range("L150").Value = 2500
$ws.Range("M150").Value = 2400
$ws.Range("P150").Value = 1600
$ws.Range("D151").Value = 44280
$ws.Range("J151").Value = 3000
$ws.Range("K151").Value = 2000
$ws.Range("L151").Value = 2500
$ws.Range("M151").Value = 2250
$ws.Range("P151").Value = 1500
$ws.Range("D152").Value = 44385
$ws.Range("J152").Value = 3320
$ws.Range("K152").Value = 1500
$ws.Range("L152").Value = 2000
$ws.Range("M152").Value = 1750
$ws.Range("P152").Value = 1167
$ws.Range("D153").Value = 44586
$ws.Range("J153").Value = 2800
$ws.Range("K153").Value = 2500
$ws.Range("L153").Value = 3000
$ws.Range("M153").Value = 2750
$ws.Range("P153").Value = 1833
$ws.Range("D154").Value = 44707
$ws.Range("J154").Value = 3000
$ws.Range("K154").Value = 1500
$ws.Range("L154").Value = 2000
$ws.Range("M154").Value = 1750
$ws.Range("P154").Value = 1167
$ws.Range("D155").Value = 44705
$ws.Range("J155").Value = 2400
$ws.Range("K155").Value = 1500
$ws.Range("L155").Value = 2000
$ws.Range("M155").Value = 1750
$ws.Range("P155").Value = 1167
$ws.Range("D156").Value = 44383
$ws.Range("J156").Value = 3200
$ws.Range("K156").Value = 1500
$ws.Range("L156").Value = 2000
$ws.Range("M156").Value = 1750
$ws.Range("P156").Value = 1167
$ws.Range("D157").Value = 44784
$ws.Range("J157").Value = 2600
$ws.Range("K157").Value = 2000
$ws.Range("L157").Value = 2500
$ws.Range("M157").Value = 2250
$ws.Range("P157").Value = 1500
$ws.Range("D158").Value = 44166
$ws.Range("J158").Value = 2800
$ws.Range("K158").Value = 1300
$ws.Range("L158").Value = 1500
$ws.Range("M158").Value = 1400
$ws.Range("P158").Value = 933
$ws.Range("D159").Value = 44476
$ws.Range("J159").Value = 3000
$ws.Range("K159").Value = 1500
$ws.Range("L159").Value = 2000
$ws.Range("M159").Value = 1750
$ws.Range("P159").Value = 1167
$ws.Range("D160").Value = 44798
$ws.Range("J160").Value = 2000
$ws.Range("K160").Value = 2000
$ws.Range("L160").Value = 2500
$ws.Range("M160").Value = 2250
$ws.Range("P160").Value = 1500
$ws.Range("D161").Value = 44579
$ws.Range("J161").Value = 3000
$ws.Range("K161").Value = 2500
$ws.Range("L161").Value = 3000
$ws.Range("M161").Value = 2750
$ws.Range("P161").Value = 1833
$ws.Range("D162").Value = 44161
$ws.Range("J162").Value = 3100
$ws.Range("K162").Value = 1300
$ws.Range("L162").Value = 1500
$ws.Range("M162").Value = 1400
$ws.Range("P162").Value = 933
$ws.Range("D163").Value = 44882
$ws.Range("J163").Value = 2600
$ws.Range("K163").Value = 1500
$ws.Range("L163").Value = 2000
$ws.Range("M163").Value = 1750
$ws.Range("P163").Value = 1167
$ws.Range("D164").Value = 44264
$ws.Range("J164").Value = 3600
$ws.Range("K164").Value = 2000
$ws.Range("L164").Value = 2500
$ws.Range("M164").Value = 2250
$ws.Range("P164").Value = 1500
$ws.Range("D165").Value = 44250
$ws.Range("J165").Value = 3400
$ws.Range("K165").Value = 1500
$ws.Range("L165").Value = 2000
$ws.Range("M165").Value = 1750
$ws.Range("P165").Value = 1167
$ws.Range("D166").Value = 44215
$ws.Range("J166").Value = 2800
$ws.Range("K166").Value = 1300
$ws.Range("L166").Value = 1500
$ws.Range("M166").Value = 1400
$ws.Range("P166").Value = 933
$ws.Range("D167").Value = 44782
$ws.Range("J167").Value = 2600
$ws.Range("K167").Value = 2000
$ws.Range("L167").Value = 2500
$ws.Range("M167").Value = 2250
$ws.Range("P167").Value = 1500
$ws.Range("D168").Value = 44754
$ws.Range("J168").Value = 2400
$ws.Range("K168").Value = 1500
$ws.Range("L168").Value = 2000
$ws.Range("M168").Value = 1750
$ws.Range("P168").Value = 1167
$ws.Range("D169").Value = 44329
$ws.Range("J169").Value = 3300
$ws.Range("K169").Value = 1300
$ws.Range("L169").Value = 1500
$ws.Range("M169").Value = 1400
$ws.Range("P169").Value = 933
$ws.Range("D170").Value = 44663
$ws.Range("J170").Value = 2360
$ws.Range("K170").Value = 2000
$ws.Range("L170").Value = 2500
$ws.Range("M170").Value = 2250
$ws.Range("P170").Value = 1500
$ws.Range("D171").Value = 44399
$ws.Range("J171").Value = 3320
$ws.Range("K171").Value = 1500
$ws.Range("L171").Value = 2000
$ws.Range("M171").Value = 1750
$ws.Range("P171").Value = 1167
$ws.Range("D172").Value = 44266
$ws.Range("J172").Value = 3600
$ws.Range("K172").Value = 2000
$ws.Range("L172").Value = 2500
$ws.Range("M172").Value = 2250
$ws.Range("P172").Value = 1500
$ws.Range("D173").Value = 44446
$ws.Range("J173").Value = 3200
$ws.Range("K173").Value = 2000
$ws.Range("L173").Value = 2500
$ws.Range("M173").Value = 2250
$ws.Range("P173").Value = 1500
$ws.Range("D174").Value = 44540
$ws.Range("J174").Value = 3000
$ws.Range("K174").Value = 1500
$ws.Range("L174").Value = 2000
$ws.Range("M174").Value = 1750
$ws.Range("P174").Value = 1167
$ws.Range("D175").Value = 44607
$ws.Range("J175").Value = 2400
$ws.Range("K175").Value = 2300
$ws.Range("L175").Value = 2500
$ws.Range("M175").Value = 2400
$ws.Range("P175").Value = 1600
$ws.Range("D176").Value = 44658
$ws.Range("J176").Value = 3000
$ws.Range("K176").Value = 2000
$ws.Range("L176").Value = 2500
$ws.Range("M176").Value = 2250
$ws.Range("P176").Value = 1500
$ws.Range("D177").Value = 44911
$ws.Range("J177").Value = 1800
$ws.Range("K177").Value = 3000
$ws.Range("L177").Value = 3500
$ws.Range("M177").Value = 3250
$ws.Range("P177").Value = 2167
$ws.Range("D178").Value = 44901
$ws.Range("J178").Value = 1600
$ws.Range("K178").Value = 3000
$ws.Range("L178").Value = 3500
$ws.Range("M178").Value = 3250
$ws.Range("P178").Value = 2167
$ws.Range("D179").Value = 44467
$ws.Range("J179").Value = 3100
$ws.Range("K179").Value = 1500
$ws.Range("L179").Value = 2000
$ws.Range("M179").Value = 1750
$ws.Range("P179").Value = 1167
$ws.Range("D180").Value = 44628
$ws.Range("J180").Value = 2400
$ws.Range("K180").Value = 2500
$ws.Range("L180").Value = 3000
$ws.Range("M180").Value = 2750
$ws.Range("P180").Value = 1833
$ws.Range("D181").Value = 44644
$ws.Range("J181").Value = 2400
$ws.Range("K181").Value = 2300
$ws.Range("L181").Value = 2500
$ws.Range("M181").Value = 2400
$ws.Range("P181").Value = 1600
$ws.Range("D182").Value = 44224
$ws.Range("J182").Value = 2800
$ws.Range("K182").Value = 1300
$ws.Range("L182").Value = 1500
$ws.Range("M182").Value = 1400
$ws.Range("P182").Value = 933
$ws.Range("D183").Value = 44637
$ws.Range("J183").Value = 2460
$ws.Range("K183").Value = 2500
$ws.Range("L183").Value = 3000
$ws.Range("M183").Value = 2750
$ws.Range("P183").Value = 1833
$ws.Range("D184").Value = 44616
$ws.Range("J184").Value = 2400
$ws.Range("K184").Value = 2300
$ws.Range("L184").Value = 2500
$ws.Range("M184").Value = 2400
$ws.Range("P184").Value = 1600
$ws.Range("D185").Value = 44581
$ws.Range("J185").Value = 3100
$ws.Range("K185").Value = 2500
$ws.Range("L185").Value = 3000
$ws.Range("M185").Value = 2750
$ws.Range("P185").Value = 1833
$ws.Range("D186").Value = 44271
$ws.Range("J186").Value = 3200
$ws.Range("K186").Value = 2000
$ws.Range("L186").Value = 2500
$ws.Range("M186").Value = 2250
$ws.Range("P186").Value = 1500
$ws.Range("D187").Value = 44908
$ws.Range("J187").Value = 2000
$ws.Range("K187").Value = 3000
$ws.Range("L187").Value = 3500
$ws.Range("M187").Value = 3250
$ws.Range("P187").Value = 2167
$ws.Range("D188").Value = 44259
$ws.Range("J188").Value = 3400
$ws.Range("K188").Value = 2000
$ws.Range("L188").Value = 2500
$ws.Range("M188").Value = 2250
$ws.Range("P188").Value = 1500
$ws.Range("D189").Value = 44252
$ws.Range("J189").Value = 3600
$ws.Range("K189").Value = 1500
$ws.Range("L189").Value = 2000
$ws.Range("M189").Value = 1750
$ws.Range("P189").Value = 1167
$ws.Range("D190").Value = 44243
$ws.Range("J190").Value = 3200
$ws.Range("K190").Value = 1500
$ws.Range("L190").Value = 2000
$ws.Range("M190").Value = 1750
$ws.Range("P190").Value = 1167
$ws.Range("D191").Value = 44539
$ws.Range("J191").Value = 3100
$ws.Range("K191").Value = 1500
$ws.Range("L191").Value = 2000
$ws.Range("M191").Value = 1750
$ws.Range("P191").Value = 1167
$ws.Range("D192").Value = 44826
$ws.Range("J192").Value = 3000
$ws.Range("K192").Value = 2000
$ws.Range("L192").Value = 2500
$ws.Range("M192").Value = 2250
$ws.Range("P192").Value = 1500
$ws.Range("D193").Value = 44838
$ws.Range("J193").Value = 2800
$ws.Range("K193").Value = 2000
$ws.Range("L193").Value = 2500
$ws.Range("M193").Value = 2250
$ws.Range("P193").Value = 1500
$ws.Range("D194").Value = 44229
$ws.Range("J194").Value = 3200
$ws.Range("K194").Value = 1300
$ws.Range("L194").Value = 1500
$ws.Range("M194").Value = 1400
$ws.Range("P194").Value = 933
$ws.Range("D195").Value = 44322
$ws.Range("J195").Value = 3320
$ws.Range("K195").Value = 1300
$ws.Range("L195").Value = 1500
$ws.Range("M195").Value = 1400
$ws.Range("P195").Value = 933
$ws.Range("D196").Value = 44320
$ws.Range("J196").Value = 3400
$ws.Range("K196").Value = 1300
$ws.Range("L196").Value = 1500
$ws.Range("M196").Value = 1400
$ws.Range("P196").Value = 933
$ws.Range("D197").Value = 44371
$ws.Range("J197").Value = 3300
$ws.Range("K197").Value = 1500
$ws.Range("L197").Value = 2000
$ws.Range("M197").Value = 1750
$ws.Range("P197").Value = 1167
$ws.Range("D198").Value = 44810
$ws.Range("J198").Value = 2400
$ws.Range("K198").Value = 2000
$ws.Range("L198").Value = 2500
$ws.Range("M198").Value = 2250
$ws.Range("P198").Value = 1500
$ws.Range("D199").Value = 44175
$ws.Range("J199").Value = 3000
$ws.Range("K199").Value = 1300
$ws.Range("L199").Value = 1500
$ws.Range("M199").Value = 1400
$ws.Range("P199").Value = 933

# New row 200 (full row, holds what used to be the last row, 199)
$ws.Range("A200").Value = 8
$ws.Range("B200").Value = "Terminal La Palmera de La Serena"
$ws.Range("C200").Value = "Coquimbo"
$ws.Range("D200").Value = 44595
$ws.Range("E200").Value = 4
$ws.Range("F200").Value = 100112044
$ws.Range("G200").Value = "Perejil"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 3000
$ws.Range("K200").Value = 2500
$ws.Range("L200").Value = 2800
$ws.Range("M200").Value = 2650
$ws.Range("N200").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O200").Value = "Provincia del Elquí"
$ws.Range("P200").Value = 1767
$ws.Range("Q200").Value = 1.5
$ws.Range("R200").Value = "Hortaliza"

# Match date formatting/style used by the rest of column D
$ws.Range("D200").NumberFormat = $ws.Range("D199").NumberFormat
